$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds plain-text price strings (European-style dot separators).
# Prefix with an apostrophe so Excel stores them as text instead of auto-
# converting number-looking values (e.g. "305.82", "1.000") to numerics,
# then reset the style back to Normal so no stray quote-prefix style sticks
# to the cell (keeps the cell at the default/unstyled format, like the original).
$ws.Range("D2").Value = "'26.961.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "'1.860.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'305.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.5058"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("D8").Value = "'0.3736"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "'0.07131"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.8874"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'20.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("D12").Value = "'0.07551"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "'1.851.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "'5.292"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "'88.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'0.000008362"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("D18").Value = "'14.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'27.002.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").Value = "'5.045"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'2.090.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "'6.462"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "'1.845"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "'147.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'2.090"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.90%  "
$ws.Range("D29").Value = "'112.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "'4.648"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.59%  "
$ws.Range("D31").Value = "'4.638"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("D32").Value = "'0.09040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").Value = "'0.05109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").Value = "'3.054"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("D35").Value = "'1.149"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.32%  "
$ws.Range("D36").Value = "'0.7281"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.74%  "
$ws.Range("D37").Value = "'0.02037"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("D38").Value = "'3.035"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  -5.97%  "
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("D41").Value = "'0.5316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").Value = "'6.584"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "'115.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("D44").Value = "'8.284"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "'0.1471"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "'0.4606"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'9.976"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.13%  "
$ws.Range("D49").Value = "'1.556"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").Value = "'36.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  -3.83%  "
